# Rotate the CLIENTE rows (113..166) on both "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets: each row's CLIENTE name plus its numeric data
# (everything in columns B..last-col) moves down by one row, and the
# values that were on the last row (166) wrap around to the first row
# (113). Column A (ASESOR) simply inherits the value that was directly
# above it, which only actually changes anything at the two rows where
# an ASESOR block boundary sits (141 and 162).

$wb = $excel.ActiveWorkbook

function Rotate-ClientBlock($ws, $lastColLetter) {
    # Read the seed row (112) together with the full block (113..166)
    # in one shot: rows 112..166 => 55 rows.
    $oldRange = $ws.Range("A112:" + $lastColLetter + "166")
    $old = $oldRange.Value()

    # Grab the destination range's own current value first so we get a
    # COM-style 1-based [,] array (matches $old's indexing) rather than
    # a 0-based .NET array.
    $targetRange = $ws.Range("A113:" + $lastColLetter + "166")
    $new = $targetRange.Value()

    $nCols = $old.GetUpperBound(1)

    for ($r = 1; $r -le 54; $r++) {
        for ($c = 1; $c -le $nCols; $c++) {
            if ($c -eq 1) {
                # Column A: value comes from the row immediately above.
                $new[$r, $c] = $old[$r, $c]
            } elseif ($r -eq 1) {
                # First row of the block wraps around to the last row
                # of the old block for every column except A.
                $new[$r, $c] = $old[55, $c]
            } else {
                $new[$r, $c] = $old[$r, $c]
            }
        }
    }

    $targetRange.Value = $new
}

$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
Rotate-ClientBlock $ws1 "R"

$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
Rotate-ClientBlock $ws2 "G"
